$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.584138035774231
$ws.Range("B1").Value = 0.4565820395946503
$ws.Range("C1").Value = 3.70927882194519
$ws.Range("D1").Value = 3.600527286529541
$ws.Range("E1").Value = 0.9671239256858826
